$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: row number, new nombre_aides (column C), new montant_total (column D)
$updates = @(
    @(2, 329058, 419105950),
    @(8, 876, 1288399),
    @(10, 119243, 174713253),
    @(12, 61080, 88144701),
    @(16, 4050, 5749637),
    @(20, 7109, 9925961),
    @(22, 79172, 98626700),
    @(28, 32965, 48247181),
    @(30, 11737, 16882730),
    @(33, 1573, 2209891),
    @(35, 1954, 2757992),
    @(36, 99154, 124657883),
    @(44, 45048, 66012119),
    @(46, 9373, 13441824),
    @(48, 1430, 1987311),
    @(51, 2481, 3468587),
    @(52, 70548, 88459159),
    @(59, 28669, 42045742),
    @(62, 11407, 16492001),
    @(68, 1579, 2213851),
    @(70, 20929, 27404826),
    @(74, 7734, 11326800),
    @(76, 5228, 7592494),
    @(78, 292, 411083),
    @(79, 144006, 179451635),
    @(83, 440, 642824),
    @(85, 64648, 94740475),
    @(88, 30393, 43967559),
    @(91, 3000, 4242243),
    @(92, 34481, 46749197),
    @(96, 8359, 12288126),
    @(98, 7729, 11216237),
    @(100, 551, 782406),
    @(101, 521, 752050),
    @(102, 11428, 18274439),
    @(104, 2792, 4797354),
    @(106, 3825, 6581934),
    @(108, 168, 285045),
    @(109, 221, 359030),
    @(110, 145011, 179340364),
    @(114, 961, 1409977),
    @(116, 53687, 78682945),
    @(118, 27947, 40490219),
    @(119, 1326, 1813694),
    @(122, 2398, 3372326),
    @(124, 538400, 711455433),
    @(129, 1409, 2088214),
    @(131, 214404, 315157910),
    @(132, 423, 631210),
    @(134, 192284, 279636318),
    @(137, 2887, 4053072),
    @(140, 6792, 9588568),
    @(143, 46095, 61516197),
    @(145, 25, 36230),
    @(149, 14409, 21118939),
    @(150, 3862, 5569092),
    @(155, 414, 583813),
    @(156, 18156, 24002421),
    @(160, 7407, 10777403),
    @(162, 5180, 7455538),
    @(167, 21046, 36877480),
    @(168, 2256, 3955071),
    @(169, 294, 501089),
    @(172, 119, 217449),
    @(173, 89714, 112051640),
    @(180, 34471, 50543665),
    @(182, 13356, 19297802),
    @(186, 1746, 2450429),
    @(188, 243637, 302682416),
    @(194, 893, 1313845),
    @(196, 88011, 128995603),
    @(199, 33795, 48656175),
    @(202, 5190, 7389358),
    @(205, 5185, 7180301),
    @(208, 270055, 334122123),
    @(215, 624, 908878),
    @(217, 96883, 141727250),
    @(220, 52882, 76425479),
    @(222, 18, 25660),
    @(223, 4736, 6647035),
    @(226, 6155, 8528238),
    @(229, 108921, 136144884),
    @(234, 571, 834439),
    @(236, 50343, 73744891),
    @(238, 12820, 18439035),
    @(242, 2662, 3729013),
    @(243, 264596, 334067782),
    @(249, 845, 1240904),
    @(251, 97722, 143185763),
    @(254, 67095, 97252397),
    @(256, 2451, 3457224),
    @(259, 4920, 6908304)
)

foreach ($u in $updates) {
    $row = $u[0]
    $newC = $u[1]
    $newD = $u[2]
    $ws.Cells.Item($row, 3).Value = $newC
    $ws.Cells.Item($row, 4).Value = $newD
}
